# "added 4wk low sales check"
# Refreshes the forecast numbers on the "Forecast Comparison" sheet (new
# MyForecast / Inventory Coverage / Seasonality Index figures, plus a
# couple of Stockout Risk / Reorder Urgency flags that flipped once the
# 4-week low-sales check was added) and rolls the updated totals through
# to the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: per-row updates ------------------------------
# Row 2 (W10)
$ws.Range("D2").Value = 70
$ws.Range("H2").Value = 5.31
$ws.Range("L2").Value = 0.8

# Row 3 (W11)
$ws.Range("D3").Value = 78
$ws.Range("H3").Value = 3.85
$ws.Range("L3").Value = 0.86

# Row 4 (W12)
$ws.Range("D4").Value = 83
$ws.Range("H4").Value = 2.7
$ws.Range("L4").Value = 0.99

# Row 5 (W13)
$ws.Range("D5").Value = 80
$ws.Range("H5").Value = 1.76
$ws.Range("L5").Value = 0.8100000000000001

# Row 6 (W14)
$ws.Range("D6").Value = 71
$ws.Range("H6").Value = 0.85
$ws.Range("J6").Value = "Urgent"
$ws.Range("L6").Value = 0.89

# Row 7 (W15)
$ws.Range("D7").Value = 65
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "High"
$ws.Range("L7").Value = 0.89

# Row 8 (W16)
$ws.Range("D8").Value = 67
$ws.Range("L8").Value = 1.04

# Row 9 (W17)
$ws.Range("D9").Value = 76
$ws.Range("L9").Value = 1.01

# Row 10 (W18)
$ws.Range("D10").Value = 83
$ws.Range("L10").Value = 1.09

# Row 11 (W19)
$ws.Range("D11").Value = 81
$ws.Range("L11").Value = 1.18

# Row 12 (W20)
$ws.Range("D12").Value = 74
$ws.Range("L12").Value = 0.96

# Row 13 (W21)
$ws.Range("D13").Value = 69
$ws.Range("L13").Value = 0.91

# Row 14 (W22)
$ws.Range("D14").Value = 74
$ws.Range("L14").Value = 1.06

# Row 15 (W23)
$ws.Range("D15").Value = 83
$ws.Range("L15").Value = 0.8

# Row 16 (W24)
$ws.Range("D16").Value = 87
$ws.Range("L16").Value = 1.18

# Row 17 (W25)
$ws.Range("D17").Value = 82
$ws.Range("L17").Value = 1.11

# --- Summary: totals/extremes recomputed from the new forecast column --
$summary.Range("B9").Value = "1227"
$summary.Range("B10").Value = "593"
$summary.Range("B11").Value = "313"
$summary.Range("B12").Value = "87"
$summary.Range("B14").Value = "65"
